# Weekly update: insert two new price records for "Macroferia Regional de
# Talca - Repollo" right before the current row 491, shifting the existing
# 491-537 data block down by two rows (to 493-539).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 491 (Excel copies formatting, incl. the
# date number format on column D, from the row immediately above).
$ws.Rows.Item(491).Resize(2).Insert()

# New row 491: Primera, 2023-08-28
$ws.Cells.Item(491, 1).Value  = 5
$ws.Cells.Item(491, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(491, 3).Value  = "Maule"
$ws.Cells.Item(491, 4).Value  = 45166
$ws.Cells.Item(491, 5).Value  = 7
$ws.Cells.Item(491, 6).Value  = 100112006
$ws.Cells.Item(491, 7).Value  = "Repollo"
$ws.Cells.Item(491, 8).Value  = "Crespo record"
$ws.Cells.Item(491, 9).Value  = "Primera"
$ws.Cells.Item(491, 10).Value = 3000
$ws.Cells.Item(491, 11).Value = 600
$ws.Cells.Item(491, 12).Value = 600
$ws.Cells.Item(491, 13).Value = 600
$ws.Cells.Item(491, 14).Value = "$/unidad"
$ws.Cells.Item(491, 15).Value = "Región del Maule"
$ws.Cells.Item(491, 16).Value = 600
$ws.Cells.Item(491, 17).Value = 1
$ws.Cells.Item(491, 18).Value = "Hortaliza"

# New row 492: Segunda, 2023-08-28
$ws.Cells.Item(492, 1).Value  = 5
$ws.Cells.Item(492, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(492, 3).Value  = "Maule"
$ws.Cells.Item(492, 4).Value  = 45166
$ws.Cells.Item(492, 5).Value  = 7
$ws.Cells.Item(492, 6).Value  = 100112006
$ws.Cells.Item(492, 7).Value  = "Repollo"
$ws.Cells.Item(492, 8).Value  = "Crespo record"
$ws.Cells.Item(492, 9).Value  = "Segunda"
$ws.Cells.Item(492, 10).Value = 3000
$ws.Cells.Item(492, 11).Value = 400
$ws.Cells.Item(492, 12).Value = 400
$ws.Cells.Item(492, 13).Value = 400
$ws.Cells.Item(492, 14).Value = "$/unidad"
$ws.Cells.Item(492, 15).Value = "Región del Maule"
$ws.Cells.Item(492, 16).Value = 400
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"
